# Add 2022-Q4 data
# 1) Insert a new "2022-Q4" worksheet (duplicate of the "2022-Q3" sheet so
#    styles/formatting match exactly) before the existing "2022-Q3" sheet,
#    then overwrite its values with the 2022-Q4 figures.
# 2) Insert a new row in the "总计" (summary) sheet for the 2022-Q4 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the "2022-Q3" sheet to use as a formatting template
# for the new "2022-Q4" sheet, placed right before it.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Row 2 (fund 010428)
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "010428"
$q4.Range("B2").ClearFormats()

$q4.Range("C2").Value = "兴银策略智选混合C"

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.37"
$q4.Range("D2").ClearFormats()

$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "92.14"
$q4.Range("E2").ClearFormats()

$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "3.67"
$q4.Range("F2").ClearFormats()

$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0136"
$q4.Range("G2").ClearFormats()

$q4.Range("H2").Value = 5

# Row 3 (fund 010427)
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "010427"
$q4.Range("B3").ClearFormats()

$q4.Range("C3").Value = "兴银策略智选混合A"

$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.24"
$q4.Range("D3").ClearFormats()

$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "92.14"
$q4.Range("E3").ClearFormats()

$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "3.67"
$q4.Range("F3").ClearFormats()

$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0088"
$q4.Range("G3").ClearFormats()

$q4.Range("H3").Value = 5

# Row 4 (fund 001730) - needs to be newly inserted
$q4.Rows.Item(4).Insert()
$q4.Range("A3").Copy($q4.Range("A4"))
$q4.Range("A4").Value = 2

$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "001730"
$q4.Range("B4").ClearFormats()

$q4.Range("C4").Value = "兴银大健康灵活配置混合"

$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "0.15"
$q4.Range("D4").ClearFormats()

$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "91.17"
$q4.Range("E4").ClearFormats()

$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "3.58"
$q4.Range("F4").ClearFormats()

$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.0054"
$q4.Range("G4").ClearFormats()

$q4.Range("H4").Value = 5

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 row into the "总计" (summary) sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("A2").Value = 0

# The row-index column (A) holds the 0-based sequence number; every row
# that got pushed down one position needs its number bumped by one.
for ($r = 3; $r -le 9; $r++) {
    $cell = $summary.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.03
$summary.Range("B2:D2").ClearFormats()

Write-Output "done"
